# Auto-generated Excel COM-interop script
# Applies scheduled-runner market price / leve profit updates across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 161.42857
$ws.Range("I5").Value = 59.75
$ws.Range("J5").Value = 297
$ws.Range("K5").Value = 59.75
$ws.Range("L5").Value = 297
$ws.Range("M5").Value = 55.25
$ws.Range("N5").Value = -527
$ws.Range("H62").Value = 2929.9
$ws.Range("I62").Value = 2787.375
$ws.Range("K62").Value = 2787.375
$ws.Range("M62").Value = -2163.375
$ws.Range("H65").Value = 2929.9
$ws.Range("I65").Value = 2787.375
$ws.Range("K65").Value = 13936.875
$ws.Range("M65").Value = -10816.875
$ws.Range("H86").Value = 2540.7693
$ws.Range("I86").Value = 2277.5334
$ws.Range("J86").Value = 2899.7273
$ws.Range("K86").Value = 2277.5334
$ws.Range("L86").Value = 2899.7273
$ws.Range("M86").Value = -1154.5334
$ws.Range("N86").Value = -5145.7273
$ws.Range("H89").Value = 2540.7693
$ws.Range("I89").Value = 2277.5334
$ws.Range("J89").Value = 2899.7273
$ws.Range("K89").Value = 11387.667
$ws.Range("L89").Value = 14498.6365
$ws.Range("M89").Value = -5771.666999999999
$ws.Range("N89").Value = -25730.6365
$ws.Range("H103").Value = 2763.125
$ws.Range("J103").Value = 2642.75
$ws.Range("L103").Value = 7928.25
$ws.Range("N103").Value = -9100.25
$ws.Range("H112").Value = 3116.5
$ws.Range("J112").Value = 3116.5
$ws.Range("L112").Value = 9349.5
$ws.Range("N112").Value = -11565.5
$ws.Range("H129").Value = 2615.6785
$ws.Range("I129").Value = 12880.75
$ws.Range("J129").Value = 904.8333
$ws.Range("K129").Value = 38642.25
$ws.Range("L129").Value = 2714.4999
$ws.Range("M129").Value = -33642.25
$ws.Range("N129").Value = -12714.4999
$ws.Range("H135").Value = 3487.2974
$ws.Range("I135").Value = 790.2222
$ws.Range("J135").Value = 6042.421
$ws.Range("K135").Value = 7111.999800000001
$ws.Range("L135").Value = 54381.789
$ws.Range("M135").Value = -4576.999800000001
$ws.Range("N135").Value = -59451.789
$ws.Range("H137").Value = 1470.6578
$ws.Range("I137").Value = 1181.5
$ws.Range("J137").Value = 2280.3
$ws.Range("K137").Value = 3544.5
$ws.Range("L137").Value = 6840.900000000001
$ws.Range("M137").Value = -994.5
$ws.Range("N137").Value = -11940.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 20866.666
$ws.Range("J3").Value = 29800
$ws.Range("L3").Value = 29800
$ws.Range("N3").Value = -30030
$ws.Range("H32").Value = 25631.51
$ws.Range("I32").Value = 8922.139999999999
$ws.Range("K32").Value = 8922.139999999999
$ws.Range("M32").Value = -8635.139999999999
$ws.Range("H61").Value = 1798.8649
$ws.Range("I61").Value = 1087.1666
$ws.Range("J61").Value = 2140.48
$ws.Range("K61").Value = 1087.1666
$ws.Range("L61").Value = 2140.48
$ws.Range("M61").Value = -875.1666
$ws.Range("N61").Value = -2564.48
$ws.Range("H74").Value = 1332.0714
$ws.Range("I74").Value = 1350.421
$ws.Range("J74").Value = 1293.3334
$ws.Range("K74").Value = 1350.421
$ws.Range("L74").Value = 1293.3334
$ws.Range("M74").Value = -476.421
$ws.Range("N74").Value = -3041.3334
$ws.Range("H77").Value = 1332.0714
$ws.Range("I77").Value = 1350.421
$ws.Range("J77").Value = 1293.3334
$ws.Range("K77").Value = 6752.105
$ws.Range("L77").Value = 6466.666999999999
$ws.Range("M77").Value = -2384.105
$ws.Range("N77").Value = -15202.667
$ws.Range("H122").Value = 2708.8
$ws.Range("I122").Value = 2262.7693
$ws.Range("K122").Value = 6788.3079
$ws.Range("M122").Value = -4338.3079
$ws.Range("H132").Value = 20108.117
$ws.Range("I132").Value = 28302.088
$ws.Range("J132").Value = 2975.2727
$ws.Range("K132").Value = 84906.264
$ws.Range("L132").Value = 8925.8181
$ws.Range("M132").Value = -82376.264
$ws.Range("N132").Value = -13985.8181
$ws.Range("H136").Value = 1798.8649
$ws.Range("I136").Value = 1087.1666
$ws.Range("J136").Value = 2140.48
$ws.Range("K136").Value = 3261.4998
$ws.Range("L136").Value = 6421.440000000001
$ws.Range("M136").Value = -711.4998000000001
$ws.Range("N136").Value = -11521.44

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3146.6086
$ws.Range("I134").Value = 3048.0833
$ws.Range("J134").Value = 3501.3
$ws.Range("K134").Value = 9144.249899999999
$ws.Range("L134").Value = 10503.9
$ws.Range("M134").Value = -6609.249899999999
$ws.Range("N134").Value = -15573.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 171.21428
$ws.Range("I7").Value = 77.22221999999999
$ws.Range("K7").Value = 77.22221999999999
$ws.Range("M7").Value = 35.77778000000001
$ws.Range("H58").Value = 2104.8635
$ws.Range("I58").Value = 1957.0625
$ws.Range("K58").Value = 1957.0625
$ws.Range("M58").Value = -1754.0625
$ws.Range("H122").Value = 989.75
$ws.Range("I122").Value = 966
$ws.Range("J122").Value = 1004
$ws.Range("K122").Value = 2898
$ws.Range("L122").Value = 3012
$ws.Range("M122").Value = -448
$ws.Range("N122").Value = -7912
$ws.Range("H132").Value = 1816.9073
$ws.Range("I132").Value = 1680.3
$ws.Range("J132").Value = 3524.5
$ws.Range("K132").Value = 5040.9
$ws.Range("L132").Value = 10573.5
$ws.Range("M132").Value = -2510.9
$ws.Range("N132").Value = -15633.5
$ws.Range("H134").Value = 1950.2
$ws.Range("I134").Value = 1024.75
$ws.Range("K134").Value = 3074.25
$ws.Range("M134").Value = -539.25
$ws.Range("H136").Value = 2104.8635
$ws.Range("I136").Value = 1957.0625
$ws.Range("K136").Value = 5871.1875
$ws.Range("M136").Value = -3321.1875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 4265.8184
$ws.Range("I56").Value = 4265.8184
$ws.Range("K56").Value = 4265.8184
$ws.Range("M56").Value = -3735.8184
$ws.Range("H60").Value = 250
$ws.Range("I60").Value = 50
$ws.Range("J60").Value = 350
$ws.Range("K60").Value = 150
$ws.Range("L60").Value = 1050
$ws.Range("M60").Value = 101
$ws.Range("N60").Value = -1552
$ws.Range("H87").Value = 6078.5
$ws.Range("I87").Value = 4771.3335
$ws.Range("J87").Value = 10000
$ws.Range("K87").Value = 14314.0005
$ws.Range("L87").Value = 30000
$ws.Range("M87").Value = -13066.0005
$ws.Range("N87").Value = -32496
$ws.Range("H90").Value = 6078.5
$ws.Range("I90").Value = 4771.3335
$ws.Range("J90").Value = 10000
$ws.Range("K90").Value = 42942.0015
$ws.Range("L90").Value = 90000
$ws.Range("M90").Value = -36702.0015
$ws.Range("N90").Value = -102480
$ws.Range("H113").Value = 871.1875
$ws.Range("I113").Value = 1189.8667
$ws.Range("J113").Value = 590
$ws.Range("K113").Value = 3569.6001
$ws.Range("L113").Value = 1770
$ws.Range("M113").Value = -1399.6001
$ws.Range("N113").Value = -6110
$ws.Range("H131").Value = 551351.0600000001
$ws.Range("J131").Value = 706465.4399999999
$ws.Range("L131").Value = 2119396.32
$ws.Range("N131").Value = -2129476.32

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2619.1667
$ws.Range("I132").Value = 1956.8182
$ws.Range("J132").Value = 3660
$ws.Range("K132").Value = 5870.4546
$ws.Range("L132").Value = 10980
$ws.Range("M132").Value = -3340.4546
$ws.Range("N132").Value = -16040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1612.3334
$ws.Range("I82").Value = 1276.7
$ws.Range("J82").Value = 2283.6
$ws.Range("K82").Value = 1276.7
$ws.Range("L82").Value = 2283.6
$ws.Range("M82").Value = -915.7
$ws.Range("N82").Value = -3005.6
$ws.Range("H85").Value = 1612.3334
$ws.Range("I85").Value = 1276.7
$ws.Range("J85").Value = 2283.6
$ws.Range("K85").Value = 1276.7
$ws.Range("L85").Value = 2283.6
$ws.Range("M85").Value = -28.70000000000005
$ws.Range("N85").Value = -4779.6
$ws.Range("H132").Value = 3864.6072
$ws.Range("I132").Value = 4887.067
$ws.Range("J132").Value = 2684.8462
$ws.Range("K132").Value = 14661.201
$ws.Range("L132").Value = 8054.5386
$ws.Range("M132").Value = -12131.201
$ws.Range("N132").Value = -13114.5386
$ws.Range("H136").Value = 2194.2
$ws.Range("I136").Value = 1584.12
$ws.Range("J136").Value = 3719.4
$ws.Range("K136").Value = 4752.36
$ws.Range("L136").Value = 11158.2
$ws.Range("M136").Value = -2202.36
$ws.Range("N136").Value = -16258.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1951.5
$ws.Range("I126").Value = 1720.0588
$ws.Range("J126").Value = 2513.5715
$ws.Range("K126").Value = 5160.1764
$ws.Range("L126").Value = 7540.7145
$ws.Range("M126").Value = -2690.1764
$ws.Range("N126").Value = -12480.7145
$ws.Range("H132").Value = 3695.4666
$ws.Range("I132").Value = 2009.878
$ws.Range("J132").Value = 7332.7896
$ws.Range("K132").Value = 6029.634
$ws.Range("L132").Value = 21998.3688
$ws.Range("M132").Value = -3499.634
$ws.Range("N132").Value = -27058.3688
$ws.Range("H136").Value = 22656.64
$ws.Range("I136").Value = 63344.188
$ws.Range("J136").Value = 5062.027
$ws.Range("K136").Value = 190032.564
$ws.Range("L136").Value = 15186.081
$ws.Range("M136").Value = -187482.564
$ws.Range("N136").Value = -20286.081
